$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / content updates (shared strings) ---

# Header: MANE -> NEDOSTACI
$ws.Range("F1").Value = "NEDOSTACI"

# Google Translate row (row 2)
$ws.Range("E2").Value = "veliki broj jezika, brzina"
$ws.Range("F2").Value = "potreban pristup internetu, nije klasični LLM, !etichal AI"

# DeepL row (row 3)
$ws.Range("E3").Value = "preciznost , kontekstualno razumijevanje, siguran u pogledu obrade podatka, omogućuje custom prevođenja "
$ws.Range("F3").Value = "potreban pristup internetu, manji broj jezika, nije podržan hrvatski , besplatna verzija: 500.000 zakova/month"

# google/flan-t5-base row (row 4)
$ws.Range("E4").Value = 'open-source, finetuniran, za primanje uputa i "praćenje tijeka misli"'
$ws.Range("F4").Value = "nije preporučeno korištenje bez prethodne procjene sigurnosti (!etichal AI)"

# facebook/m2m100_1.2B row (row 5)
$ws.Range("F5").Value = "manja preciznost u pojedinim jezicima, trenutno se ne može koristiti preko API-a nego samo lokalno"

# facebook/nllb-200-3.3B row (row 6)
$ws.Range("F6").Value = "smanjena kvaliteta kod dužih tekstova, ne može se koristiti za službene prijevode, nije u potpunosti spreman za produkcijsku upotrebu"

# DeepSeek V3 Base row (row 9)
$ws.Range("E9").Value = "jeftiniji od ostalih modela, pre-trained i post-trained, podržava više jezika u odnosu na chat"

# DeepSeek V3 Chat row (row 10)
$ws.Range("E10").Value = "jeftiniji od ostalih modela, fine-tuned, optimiziran za dijalog i interakciju"

# Claude 3.5 Sonnet row (row 11)
$ws.Range("E11").Value = "preciznost, razumijevanje konteksta, složeni zadaci, etichal AI, može se koristiti za profesionalno prevođenje, podržava veći broj znakova od haiku"

# Claude 3.5 Haiku row (row 12)
$ws.Range("F12").Value = "nije namijenjen prevođenju, potreban pristup internetu"

# MERGE row (row 14) - collapse rich text run into plain text "MERGE*"
$ws.Range("A14").Value = "MERGE*"

# --- View changes ---
$ws.Range("B19").Select()

# --- Column width change ---
$ws.Columns.Item(5).ColumnWidth = 125.43
